$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.396.37"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "2.327.44"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'511.79"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "'132.32"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("D11").Value = "'5.25"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "'0.337"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.745.03"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'23.51"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "56.356.06"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "2.316.17"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "'10.37"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "'323.11"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").Value = "'4.15"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").Value = "'6.60"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'61.26"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("D24").Value = "'8.55"
$ws.Range("E24").Value = "  +11.04%  "
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +6.97%  "
$ws.Range("D28").Value = "'167.08"
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").Value = "0.0₃0718"
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("D30").Value = "'1.66"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").Value = "'6.09"
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("D32").Value = "'18.27"
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("E35").Value = "  +2.29%  "
$ws.Range("D36").Value = "'3.94"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").Value = "'0.880"
$ws.Range("E37").Value = "  -5.62%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'38.39"
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.55"
$ws.Range("E39").Value = "  +2.41%  "
$ws.Range("D40").Value = "'150.23"
$ws.Range("E40").Value = "  +9.23%  "
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("D42").Value = "'3.55"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "'276.84"
$ws.Range("E43").Value = "  +1.24%  "
$ws.Range("D44").Value = "'5.04"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").Value = "'0.0495"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "'0.553"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "'18.04"
$ws.Range("E48").Value = "  +6.02%  "
$ws.Range("D49").Value = "'0.377"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "'0.0213"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("D51").Value = "'16.92"
$ws.Range("E51").Value = "  +1.56%  "
